# Strava data dictionary - insert a "merge key" row + column, add
# instructions rows, and colour-code Strava/Form/Merger provenance.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Insert a new row 2 (pushes the whole table, incl. the isolated
#    D21 style-only cell, down by one row -> D22).
# ---------------------------------------------------------------
$ws.Rows(2).Insert()

# ---------------------------------------------------------------
# 2. New header-ish row 2 content ("Key" legend row).
# ---------------------------------------------------------------
$ws.Range("B2").Value = "Key"
$ws.Range("C2").Value = "y or x"
$ws.Range("D2").Value = "x = strava, y = form. Merge if one is x and one is y"
$ws.Range("E2").Value = "Key"

# Row 1 header for column E changes wording.
$ws.Range("E1").Value = "Keep/preferred"

# ---------------------------------------------------------------
# 3. Fill in the new column E ("Strava"/"Form"/"Merger" tag) for
#    every data row (rows 3-17 after the insert).
# ---------------------------------------------------------------
$ws.Range("E3").Value = "Merger"
$ws.Range("E4").Value = "Merger"
$ws.Range("E5").Value = "Strava"
$ws.Range("E6").Value = "Strava"
$ws.Range("E7").Value = "Form"
$ws.Range("E8").Value = "Strava"
$ws.Range("E9").Value = "Strava"
$ws.Range("E10").Value = "Strava"
$ws.Range("E11").Value = "Strava"
$ws.Range("E12").Value = "Strava"
$ws.Range("E13").Value = "Strava"
$ws.Range("E14").Value = "Strava"
$ws.Range("E15").Value = "Strava"
$ws.Range("E16").Value = "Strava"
$ws.Range("E17").Value = "Form"

# ---------------------------------------------------------------
# 4. New instruction rows 19-20.
# ---------------------------------------------------------------
$ws.Range("A19").Value = "How to merge these 2 tables?"
$ws.Range("A20").Value = "If date = date, type = type and one is x and one is y [elapsed time similar, name similar, time of day similar- not reliable] then rows combine to one"

# ---------------------------------------------------------------
# Note: the old D21 style-only cell is already carried down to D22
# automatically by the row-insert shift above, so nothing else is
# needed there.
# ---------------------------------------------------------------

# ---------------------------------------------------------------
# 5. Selection moves to C22, matching the saved sheetView.
# ---------------------------------------------------------------
$ws.Range("C22").Select()
